# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Bahamut profits workbook
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 52937.047
$ws.Range("I62").Value = 81417.62
$ws.Range("J62").Value = 6656.125
$ws.Range("K62").Value = 81417.62
$ws.Range("L62").Value = 6656.125
$ws.Range("M62").Value = -80793.62
$ws.Range("N62").Value = -7904.125

$ws.Range("H65").Value = 52937.047
$ws.Range("I65").Value = 81417.62
$ws.Range("J65").Value = 6656.125
$ws.Range("K65").Value = 407088.1
$ws.Range("L65").Value = 33280.625
$ws.Range("M65").Value = -403968.1
$ws.Range("N65").Value = -39520.625

$ws.Range("H111").Value = 1016.25
$ws.Range("I111").Value = 1018.5714
$ws.Range("K111").Value = 3055.7142
$ws.Range("M111").Value = 11.28579999999965

$ws.Range("H127").Value = 2482.5
$ws.Range("I127").Value = 1248.75
$ws.Range("J127").Value = 2790.9375
$ws.Range("K127").Value = 3746.25
$ws.Range("L127").Value = 8372.8125
$ws.Range("M127").Value = 1213.75
$ws.Range("N127").Value = -18292.8125

$ws.Range("H129").Value = 824007.8
$ws.Range("I129").Value = 371.6
$ws.Range("K129").Value = 1114.8
$ws.Range("M129").Value = 3885.2

$ws.Range("H132").Value = 1853460.9
$ws.Range("I132").Value = 1399.7333
$ws.Range("J132").Value = 11113767
$ws.Range("K132").Value = 4199.199900000001
$ws.Range("L132").Value = 33341301
$ws.Range("M132").Value = -1669.199900000001
$ws.Range("N132").Value = -33346361

$ws.Range("H135").Value = 773.04254
$ws.Range("I135").Value = 756.75
$ws.Range("J135").Value = 1012
$ws.Range("K135").Value = 6810.75
$ws.Range("L135").Value = 9108
$ws.Range("M135").Value = -4275.75
$ws.Range("N135").Value = -14178

$ws.Range("H137").Value = 746.8679
$ws.Range("I137").Value = 676.2273
$ws.Range("K137").Value = 2028.6819
$ws.Range("M137").Value = 521.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6634.46
$ws.Range("I32").Value = 6463.735
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 6463.735
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -6176.735
$ws.Range("N32").Value = -15574

$ws.Range("H61").Value = 1847.102
$ws.Range("I61").Value = 1797.2858
$ws.Range("J61").Value = 1971.6428
$ws.Range("K61").Value = 1797.2858
$ws.Range("L61").Value = 1971.6428
$ws.Range("M61").Value = -1585.2858
$ws.Range("N61").Value = -2395.6428

$ws.Range("H74").Value = 1075.8049
$ws.Range("I74").Value = 1087.6428
$ws.Range("K74").Value = 1087.6428
$ws.Range("M74").Value = -213.6428000000001

$ws.Range("H77").Value = 1075.8049
$ws.Range("I77").Value = 1087.6428
$ws.Range("K77").Value = 5438.214
$ws.Range("M77").Value = -1070.214

$ws.Range("H122").Value = 655.1667
$ws.Range("I122").Value = 655.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1965.5001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 484.4999
$ws.Range("N122").ClearContents() | Out-Null

$ws.Range("H123").Value = 11785.643
$ws.Range("J123").Value = 11785.643
$ws.Range("L123").Value = 11785.643
$ws.Range("N123").Value = -21585.643

$ws.Range("H132").Value = 2091.5789
$ws.Range("I132").Value = 1689.9584
$ws.Range("J132").Value = 2780.0715
$ws.Range("K132").Value = 5069.8752
$ws.Range("L132").Value = 8340.2145
$ws.Range("M132").Value = -2539.8752
$ws.Range("N132").Value = -13400.2145

$ws.Range("H136").Value = 1847.102
$ws.Range("I136").Value = 1797.2858
$ws.Range("J136").Value = 1971.6428
$ws.Range("K136").Value = 5391.857400000001
$ws.Range("L136").Value = 5914.928400000001
$ws.Range("M136").Value = -2841.857400000001
$ws.Range("N136").Value = -11014.9284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22433.117
$ws.Range("I31").Value = 2646.8
$ws.Range("J31").Value = 65715.69
$ws.Range("K31").Value = 2646.8
$ws.Range("L31").Value = 65715.69
$ws.Range("M31").Value = -2351.8
$ws.Range("N31").Value = -66305.69

$ws.Range("H34").Value = 22433.117
$ws.Range("I34").Value = 2646.8
$ws.Range("J34").Value = 65715.69
$ws.Range("K34").Value = 2646.8
$ws.Range("L34").Value = 65715.69
$ws.Range("M34").Value = -2444.8
$ws.Range("N34").Value = -66119.69

$ws.Range("H58").Value = 5531.7036
$ws.Range("J58").Value = 10567.637
$ws.Range("L58").Value = 10567.637
$ws.Range("N58").Value = -10973.637

$ws.Range("H122").Value = 1597.1428
$ws.Range("I122").Value = 1453.2
$ws.Range("J122").Value = 1957
$ws.Range("K122").Value = 4359.6
$ws.Range("L122").Value = 5871
$ws.Range("M122").Value = -1909.6
$ws.Range("N122").Value = -10771

$ws.Range("H132").Value = 1347.1666
$ws.Range("I132").Value = 1035.0476
$ws.Range("J132").Value = 3532
$ws.Range("K132").Value = 3105.142800000001
$ws.Range("L132").Value = 10596
$ws.Range("M132").Value = -575.1428000000005
$ws.Range("N132").Value = -15656

$ws.Range("H134").Value = 31251136
$ws.Range("I134").Value = 1264.8334
$ws.Range("J134").Value = 125000750
$ws.Range("K134").Value = 3794.5002
$ws.Range("L134").Value = 375002250
$ws.Range("M134").Value = -1259.5002
$ws.Range("N134").Value = -375007320

$ws.Range("H136").Value = 5531.7036
$ws.Range("J136").Value = 10567.637
$ws.Range("L136").Value = 31702.911
$ws.Range("N136").Value = -36802.911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2963
$ws.Range("I98").Value = 150
$ws.Range("J98").Value = 3364.8572
$ws.Range("K98").Value = 450
$ws.Range("L98").Value = 10094.5716
$ws.Range("M98").Value = 1048
$ws.Range("N98").Value = -13090.5716

$ws.Range("H131").Value = 18593932
$ws.Range("J131").Value = 47042.863
$ws.Range("L131").Value = 141128.589
$ws.Range("N131").Value = -151208.589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3916.5833
$ws.Range("I70").Value = 3783.9333
$ws.Range("K70").Value = 3783.9333
$ws.Range("M70").Value = -3513.9333

$ws.Range("H73").Value = 3916.5833
$ws.Range("I73").Value = 3783.9333
$ws.Range("K73").Value = 3783.9333
$ws.Range("M73").Value = -2847.9333

$ws.Range("H122").Value = 2193848.8
$ws.Range("I122").Value = 2632418.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7897255.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -7894805.5
$ws.Range("N122").Value = -7900

$ws.Range("H132").Value = 4202.871
$ws.Range("I132").Value = 4278.3335
$ws.Range("J132").Value = 4044.4
$ws.Range("K132").Value = 12835.0005
$ws.Range("L132").Value = 12133.2
$ws.Range("M132").Value = -10305.0005
$ws.Range("N132").Value = -17193.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3257
$ws.Range("I16").Value = 3257
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3257
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3087
$ws.Range("N16").ClearContents() | Out-Null

$ws.Range("H122").Value = 2483.4285
$ws.Range("I122").Value = 2480.6667
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7442.000100000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4992.000100000001
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1719.8788
$ws.Range("I132").Value = 1529.16
$ws.Range("J132").Value = 2315.875
$ws.Range("K132").Value = 4587.48
$ws.Range("L132").Value = 6947.625
$ws.Range("M132").Value = -2057.48
$ws.Range("N132").Value = -12007.625
